# Add a new "2080" PAGE perturbation input option, mirroring the existing
# 2060 block (columns AQ:AV) into a new block in columns AX:BC, for both the
# CH4 table (rows 1-13) and the N2O table (rows 15-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- CH4 table (rows 3-13 header+data) ----
$ws.Range("AQ3:AV13").Copy()
$ws.Range("AX3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("AQ3:AV13").Copy()
$ws.Range("AX3").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

# New label above the CH4 2080 block (row 2)
$ws.Range("AX2").Value2 = "PAGE input: CH4 Shock, 2080"

# ---- N2O table (rows 17-27 header+data) ----
$ws.Range("AQ17:AV27").Copy()
$ws.Range("AX17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AQ17:AV27").Copy()
$ws.Range("AX17").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# New label above the N2O 2080 block (row 16)
$ws.Range("AX16").Value2 = "PAGE input: N2O Shock, 2080"
